# Adds new Rules section "r6.0 How to Win" (with sub-rules r6.1-r6.3) to the
# Rules sheet as four new rows inserted right after the existing row 37 (r4.95),
# pushing the rest of the table down by four rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("38:41").Insert()

# Column A: short rule-number labels
$ws.Range("A38").Value2 = "r6.0"
$ws.Range("B38").Value2 = "<Bold>r6.0 How to Win</Bold>`r`n<LineBreak/><LineBreak/>`r`nIn order to win an engagement or scenario, both you and your tank must survive. If you are killed or wounded sufficently to miss combat for any time, or if you tank is knocked out, you lose. `r`n<LineBreak/><LineBreak/>`r`nIf both you and your tank survive, winning depends on your score of victory points. Vicotry in the campaign game is a special case.`r`n<LineBreak/><LineBreak/>`r`n<InlineUIContainer><Button Content='r6.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Victory Points<LineBreak/>`r`n<InlineUIContainer><Button Content='r6.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Victory<LineBreak/>`r`n<InlineUIContainer><Button Content='r6.3' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> Campaign Game Victory"
$ws.Range("A39").Value2 = "r6.1"
$ws.Range("A40").Value2 = "r6.2"
$ws.Range("A41").Value2 = "r6.3"

# Column B: the rule text, filled in out of row order (r6.2 before r6.1, then r6.3)
$ws.Range("B40").Value2 = "<Bold>r6.2 Victory</Bold>`r`n<LineBreak/><LineBreak/>`r`nWhen the engagement or scenario ends, multiply the victory points in the four categories by the multipliers shown on the After Action Report `r`n<InlineUIContainer><Button Content='AAR' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>. `r`n<LineBreak/><LineBreak/>`r`nDepending on whether the 4th Armored Division is advancing, engaged in a pitched battle, or defending against a counterattack; the gain or loss of territory and the destruction of enemy units will be more or less important.`r`n<LineBreak/><LineBreak/>`r`nIf the sum of the adjusted positive and negative points is positive, you win. If the sum is negative, you lose."
$ws.Range("B39").Value2 = "<Bold>r6.1 Victory Points</Bold>`r`n<LineBreak/><LineBreak/>`r`nVictory points are scored for knocking out enemy units and capturing enemy territory. They are subtracted for loss of friendly forces and the loss of US controlled territory.`r`n<LineBreak/><LineBreak/>`r`nThe points scored for each acheivement are shown on the After Action Report `r`n<InlineUIContainer><Button Content='AAR' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.`r`n Each time a victory point action occurs, the value is added on the AAR. `r`n<LineBreak/><LineBreak/>`r`nWhen the scenario ends, total the victory points into four categories:<LineBreak/>`r`n-- Enemy units KO'd by your tank<LineBreak/>`r`n-- Enemy units KO'd by friendly forces<LineBreak/>`r`n-- Territory controlled or lost<LineBreak/>`r`n-- Friendly force losses (tanks and infantry squads)"
$ws.Range("B41").Value2 = "<Bold>r6.3 Campaign Game Victory</Bold>`r`n<LineBreak/><LineBreak/>`r`nThe campaign in the European theatre of operations ended in victory for the Allied armies, and the 4th Armored Division must be considered one of the big winners in view of their tremendous war record.`r`n<LineBreak/><LineBreak/>`r`nFor you as a tank commander, campaign victory is measured at a more basic level. If you are killed, you lose. Instead of starting over, it is recommended that you assume the identify of your replacement and continue from there.`r`n<LineBreak/><LineBreak/>`r`nIf you survive the entire campaign regardless of wounds, you win!"

# Row heights (auto-fit values from the source document)
$ws.Rows(38).RowHeight = 135
$ws.Rows(39).RowHeight = 195
$ws.Rows(40).RowHeight = 120
$ws.Rows(41).RowHeight = 105

# The sheet's remembered last-sort range (for the still-empty reserved rows further
# down the table) needs to shift down by the same four rows as everything else.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A140:A315"))
$ws.Sort.SetRange($ws.Range("A140:B315"))
$ws.Sort.Apply()

# Selection state matching the edited workbook
[void]$ws.Range("B39").Select()

Write-Host "Inserted r6.0-r6.3 rows"
